# Generate Report for Handoff
# Update the "Status" for the d7367dae... file from "Handed back: in sync with en-US"
# to "Ready for handoff" across the Overview, zh-cn and de-de sheets, and update the
# corresponding "Latest Handoff Datetime" timestamps on the language sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-02-18 03:13:54"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-02-18 03:14:06"
